$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add "Have all connectors for PCB" style markers in column H
$ws.Range("H9").Value = "^yes (Under RJ11)"
$ws.Range("H12").Value = "Yes"
$ws.Range("H13").Value = "Yes"
$ws.Range("H14").Value = "yes"
$ws.Range("H15").Value = "Yes"
$ws.Range("H21").Value = "Yes"
$ws.Range("H28").Value = "Yes"
$ws.Range("H30").Value = "Yes"
$ws.Range("H31").Value = "Yes"
$ws.Range("H33").Value = "Yes"
$ws.Range("H35").Value = "Yes"
$ws.Range("H36").Value = "Yes"
$ws.Range("H37").Value = "Yes"
$ws.Range("H40").Value = "Yes"
$ws.Range("H48").Value = "Yes"

# Update view position / selection to match final saved state
$ws.Application.ActiveWindow.ScrollRow = 25
$ws.Application.ActiveWindow.ScrollColumn = 1
$ws.Range("I35").Select()
